# Statusupdate KW23/24 ohne GM
#
# Slide 5 ("Python Skript für Sprachaufnahme") has a content placeholder
# ("Inhaltsplatzhalter 3") whose text body ends with two identical, empty
# trailing paragraphs. Remove the superfluous last one, leaving a single
# trailing empty paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item("Inhaltsplatzhalter 3")

$tr = $sh.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count

$lastPara = $tr.Paragraphs($paraCount, 1)
$lastPara.Delete()
